$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.567.32"
$ws.Range("E2").Value = "  -1.56%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.533.17"
$ws.Range("E3").Value = "  -1.40%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "602.28"
$ws.Range("E5").Value = "  -1.17%  "

# Row 6 - Solana
$ws.Range("D6").Value = "143.29"
$ws.Range("E6").Value = "  -2.68%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.532.31"

# Row 9 - XRP
$ws.Range("D9").Value = "0.518"
$ws.Range("E9").Value = "  +5.70%  "

# Row 10 and 11 - Toncoin/Dogecoin swapped positions
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.133"
$ws.Range("E10").Value = "  -2.34%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "7.86"
$ws.Range("E11").Value = "  -1.76%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -2.51%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.140.69"
$ws.Range("E13").Value = "  -1.17%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -5.85%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "28.44"
$ws.Range("E15").Value = "  -5.34%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.534.94"
$ws.Range("E16").Value = "  -0.78%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  +1.51%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "65.607.93"
$ws.Range("E18").Value = "  -1.57%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "11.00"
$ws.Range("E19").Value = "  -3.39%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  -1.02%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "14.41"
$ws.Range("E21").Value = "  -3.98%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "418.98"
$ws.Range("E22").Value = "  -3.02%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  -3.96%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "77.44"
$ws.Range("E24").Value = "  -2.10%  "

# Row 25 - WrappedeETH
$ws.Range("D25").Value = "3.678.48"
$ws.Range("E25").Value = "  -1.24%  "

# Row 27 - PEPE
$ws.Range("D27").Value = "0.0000114"
$ws.Range("E27").Value = "  -5.53%  "

# Row 28 - PancakeSwap
$ws.Range("E28").Value = "  -2.40%  "

# Row 29 - RenderToken
$ws.Range("D29").Value = "7.78"
$ws.Range("E29").Value = "  -4.29%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "8.84"
$ws.Range("E30").Value = "  -4.94%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.23%  "

# Row 32 - RenzoRestakedETH
$ws.Range("D32").Value = "3.546.11"
$ws.Range("E32").Value = "  -0.90%  "

# Row 33 - Kaspa
$ws.Range("D33").Value = "0.155"
$ws.Range("E33").Value = "  -0.88%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "24.21"
$ws.Range("E34").Value = "  -5.10%  "

# Row 35 - USDe
$ws.Range("E35").Value = "  -0.02%  "

# Row 36 - Fetch.AI
$ws.Range("E36").Value = "  -8.72%  "

# Row 37 - Aptos
$ws.Range("D37").Value = "7.55"
$ws.Range("E37").Value = "  -3.66%  "

# Row 38 and 39 - Monero/NEARProtocol swapped positions
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "5.23"
$ws.Range("E38").Value = "  -6.86%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "171.66"
$ws.Range("E39").Value = "  -1.29%  "

# Row 40 - ImmutableX
$ws.Range("D40").Value = "1.56"
$ws.Range("E40").Value = "  -8.50%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "0.0811"
$ws.Range("E41").Value = "  -5.07%  "

# Row 42 - Filecoin
$ws.Range("D42").Value = "5.04"
$ws.Range("E42").Value = "  -3.57%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  -4.30%  "

# Row 44 - OKB
$ws.Range("D44").Value = "45.25"
$ws.Range("E44").Value = "  -1.93%  "

# Row 45 - Stacks
$ws.Range("E45").Value = "  -7.59%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  +0.05%  "

# Row 47 and 48 - dogwifhat/EnergySwap swapped positions
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "23.49"
$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  -8.05%  "

# Row 49 - Cosmos
$ws.Range("D49").Value = "7.02"
$ws.Range("E49").Value = "  -2.58%  "

# Row 50 - ONDO
$ws.Range("D50").Value = "1.10"
$ws.Range("E50").Value = "  -7.63%  "

# Row 51 - SuiNetwork
$ws.Range("D51").Value = "0.901"
$ws.Range("E51").Value = "  -4.77%  "

$wb.Save()
